$d = $word.ActiveDocument

$d.Content.Find.Execute("46-32=14", $true, $false, $false, $false, $false, $true, 1, $false, "36+39=75", 2) | Out-Null
$d.Content.Find.Execute("19-9=10", $true, $false, $false, $false, $false, $true, 1, $false, "52-30=22", 2) | Out-Null
$d.Content.Find.Execute("62-31=31", $true, $false, $false, $false, $false, $true, 1, $false, "78+16=94", 2) | Out-Null
$d.Content.Find.Execute("4+76=80", $true, $false, $false, $false, $false, $true, 1, $false, "50-1=49", 2) | Out-Null
$d.Content.Find.Execute("46+7=53", $true, $false, $false, $false, $false, $true, 1, $false, "11+18=29", 2) | Out-Null
$d.Content.Find.Execute("81-22=59", $true, $false, $false, $false, $false, $true, 1, $false, "2+42=44", 2) | Out-Null
$d.Content.Find.Execute("98+1=99", $true, $false, $false, $false, $false, $true, 1, $false, "4+59=63", 2) | Out-Null
$d.Content.Find.Execute("46+49=95", $true, $false, $false, $false, $false, $true, 1, $false, "26+7=33", 2) | Out-Null
$d.Content.Find.Execute("44+36=80", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=24", 2) | Out-Null
$d.Content.Find.Execute("27+5=32", $true, $false, $false, $false, $false, $true, 1, $false, "34+44=78", 2) | Out-Null
$d.Content.Find.Execute("35+4=39", $true, $false, $false, $false, $false, $true, 1, $false, "50+42=92", 2) | Out-Null
$d.Content.Find.Execute("3+77=80", $true, $false, $false, $false, $false, $true, 1, $false, "51-30=21", 2) | Out-Null
$d.Content.Find.Execute("14+80=94", $true, $false, $false, $false, $false, $true, 1, $false, "57-34=23", 2) | Out-Null
$d.Content.Find.Execute("69-13=56", $true, $false, $false, $false, $false, $true, 1, $false, "80-26=54", 2) | Out-Null
$d.Content.Find.Execute("23-2=21", $true, $false, $false, $false, $false, $true, 1, $false, "44-38=6", 2) | Out-Null
$d.Content.Find.Execute("35+55=90", $true, $false, $false, $false, $false, $true, 1, $false, "59-47=12", 2) | Out-Null
$d.Content.Find.Execute("94-45=49", $true, $false, $false, $false, $false, $true, 1, $false, "10+62=72", 2) | Out-Null
$d.Content.Find.Execute("30+50=80", $true, $false, $false, $false, $false, $true, 1, $false, "85-44=41", 2) | Out-Null
$d.Content.Find.Execute("6+57=63", $true, $false, $false, $false, $false, $true, 1, $false, "5+53=58", 2) | Out-Null
$d.Content.Find.Execute("5+60=65", $true, $false, $false, $false, $false, $true, 1, $false, "11+63=74", 2) | Out-Null
$d.Content.Find.Execute("39-6=33", $true, $false, $false, $false, $false, $true, 1, $false, "65+30=95", 2) | Out-Null
$d.Content.Find.Execute("32+64=96", $true, $false, $false, $false, $false, $true, 1, $false, "83-75=8", 2) | Out-Null
$d.Content.Find.Execute("24-2=22", $true, $false, $false, $false, $false, $true, 1, $false, "50-38=12", 2) | Out-Null
$d.Content.Find.Execute("76-67=9", $true, $false, $false, $false, $false, $true, 1, $false, "6+89=95", 2) | Out-Null
$d.Content.Find.Execute("1+17=18", $true, $false, $false, $false, $false, $true, 1, $false, "63+32=95", 2) | Out-Null
$d.Content.Find.Execute("55+38=93", $true, $false, $false, $false, $false, $true, 1, $false, "49-38=11", 2) | Out-Null
$d.Content.Find.Execute("69-3=66", $true, $false, $false, $false, $false, $true, 1, $false, "40-24=16", 2) | Out-Null
$d.Content.Find.Execute("70-51=19", $true, $false, $false, $false, $false, $true, 1, $false, "6+55=61", 2) | Out-Null
$d.Content.Find.Execute("9+77=86", $true, $false, $false, $false, $false, $true, 1, $false, "3+53=56", 2) | Out-Null
$d.Content.Find.Execute("24+25=49", $true, $false, $false, $false, $false, $true, 1, $false, "86-81=5", 2) | Out-Null
$d.Content.Find.Execute("86-17=69", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=21", 2) | Out-Null
$d.Content.Find.Execute("44+12=56", $true, $false, $false, $false, $false, $true, 1, $false, "28+58=86", 2) | Out-Null
$d.Content.Find.Execute("42-32=10", $true, $false, $false, $false, $false, $true, 1, $false, "41-22=19", 2) | Out-Null
$d.Content.Find.Execute("32+60=92", $true, $false, $false, $false, $false, $true, 1, $false, "31+56=87", 2) | Out-Null
$d.Content.Find.Execute("4+71=75", $true, $false, $false, $false, $false, $true, 1, $false, "35+48=83", 2) | Out-Null
$d.Content.Find.Execute("56+7=63", $true, $false, $false, $false, $false, $true, 1, $false, "15-12=3", 2) | Out-Null
$d.Content.Find.Execute("32+13=45", $true, $false, $false, $false, $false, $true, 1, $false, "45+22=67", 2) | Out-Null
$d.Content.Find.Execute("90-2=88", $true, $false, $false, $false, $false, $true, 1, $false, "13+74=87", 2) | Out-Null
$d.Content.Find.Execute("81-43=38", $true, $false, $false, $false, $false, $true, 1, $false, "98-77=21", 2) | Out-Null
$d.Content.Find.Execute("91-66=25", $true, $false, $false, $false, $false, $true, 1, $false, "66+10=76", 2) | Out-Null
$d.Content.Find.Execute("34+47=81", $true, $false, $false, $false, $false, $true, 1, $false, "46-38=8", 2) | Out-Null
$d.Content.Find.Execute("78-70=8", $true, $false, $false, $false, $false, $true, 1, $false, "57-7=50", 2) | Out-Null
$d.Content.Find.Execute("98-82=16", $true, $false, $false, $false, $false, $true, 1, $false, "89-85=4", 2) | Out-Null
$d.Content.Find.Execute("96-43=53", $true, $false, $false, $false, $false, $true, 1, $false, "65-43=22", 2) | Out-Null
$d.Content.Find.Execute("51-24=27", $true, $false, $false, $false, $false, $true, 1, $false, "51-42=9", 2) | Out-Null
$d.Content.Find.Execute("17+22=39", $true, $false, $false, $false, $false, $true, 1, $false, "53+19=72", 2) | Out-Null
$d.Content.Find.Execute("28+46=74", $true, $false, $false, $false, $false, $true, 1, $false, "0+72=72", 2) | Out-Null
$d.Content.Find.Execute("72-57=15", $true, $false, $false, $false, $false, $true, 1, $false, "31-26=5", 2) | Out-Null
$d.Content.Find.Execute("21+55=76", $true, $false, $false, $false, $false, $true, 1, $false, "87-13=74", 2) | Out-Null
$d.Content.Find.Execute("73-67=6", $true, $false, $false, $false, $false, $true, 1, $false, "8+33=41", 2) | Out-Null
$d.Content.Find.Execute("46+40=86", $true, $false, $false, $false, $false, $true, 1, $false, "53+43=96", 2) | Out-Null
$d.Content.Find.Execute("95-78=17", $true, $false, $false, $false, $false, $true, 1, $false, "6+68=74", 2) | Out-Null
$d.Content.Find.Execute("7+50=57", $true, $false, $false, $false, $false, $true, 1, $false, "40+56=96", 2) | Out-Null
$d.Content.Find.Execute("64+5=69", $true, $false, $false, $false, $false, $true, 1, $false, "60-40=20", 2) | Out-Null
$d.Content.Find.Execute("37+15=52", $true, $false, $false, $false, $false, $true, 1, $false, "71-25=46", 2) | Out-Null
$d.Content.Find.Execute("8+65=73", $true, $false, $false, $false, $false, $true, 1, $false, "37+39=76", 2) | Out-Null
$d.Content.Find.Execute("68+4=72", $true, $false, $false, $false, $false, $true, 1, $false, "38+43=81", 2) | Out-Null
$d.Content.Find.Execute("18+61=79", $true, $false, $false, $false, $false, $true, 1, $false, "80-54=26", 2) | Out-Null
$d.Content.Find.Execute("85-33=52", $true, $false, $false, $false, $false, $true, 1, $false, "81-55=26", 2) | Out-Null
$d.Content.Find.Execute("25+65=90", $true, $false, $false, $false, $false, $true, 1, $false, "51+16=67", 2) | Out-Null
$d.Content.Find.Execute("19-1=18", $true, $false, $false, $false, $false, $true, 1, $false, "29+53=82", 2) | Out-Null
$d.Content.Find.Execute("85-74=11", $true, $false, $false, $false, $false, $true, 1, $false, "34-14=20", 2) | Out-Null
$d.Content.Find.Execute("14+31=45", $true, $false, $false, $false, $false, $true, 1, $false, "86-35=51", 2) | Out-Null
$d.Content.Find.Execute("20+46=66", $true, $false, $false, $false, $false, $true, 1, $false, "80-10=70", 2) | Out-Null
$d.Content.Find.Execute("99-79=20", $true, $false, $false, $false, $false, $true, 1, $false, "2+53=55", 2) | Out-Null
$d.Content.Find.Execute("28+29=57", $true, $false, $false, $false, $false, $true, 1, $false, "54+15=69", 2) | Out-Null
$d.Content.Find.Execute("46+3=49", $true, $false, $false, $false, $false, $true, 1, $false, "75-45=30", 2) | Out-Null
$d.Content.Find.Execute("18+8=26", $true, $false, $false, $false, $false, $true, 1, $false, "98-13=85", 2) | Out-Null
$d.Content.Find.Execute("42+27=69", $true, $false, $false, $false, $false, $true, 1, $false, "73-12=61", 2) | Out-Null
$d.Content.Find.Execute("76-19=57", $true, $false, $false, $false, $false, $true, 1, $false, "98-62=36", 2) | Out-Null
$d.Content.Find.Execute("75-62=13", $true, $false, $false, $false, $false, $true, 1, $false, "68+18=86", 2) | Out-Null
$d.Content.Find.Execute("52-0=52", $true, $false, $false, $false, $false, $true, 1, $false, "49-7=42", 2) | Out-Null
$d.Content.Find.Execute("46-44=2", $true, $false, $false, $false, $false, $true, 1, $false, "54+37=91", 2) | Out-Null
$d.Content.Find.Execute("65-3=62", $true, $false, $false, $false, $false, $true, 1, $false, "84-33=51", 2) | Out-Null
$d.Content.Find.Execute("21+64=85", $true, $false, $false, $false, $false, $true, 1, $false, "11+48=59", 2) | Out-Null
$d.Content.Find.Execute("65-23=42", $true, $false, $false, $false, $false, $true, 1, $false, "50+42=92", 2) | Out-Null
$d.Content.Find.Execute("95-13=82", $true, $false, $false, $false, $false, $true, 1, $false, "70-9=61", 2) | Out-Null
$d.Content.Find.Execute("57-45=12", $true, $false, $false, $false, $false, $true, 1, $false, "81-80=1", 2) | Out-Null
$d.Content.Find.Execute("45+38=83", $true, $false, $false, $false, $false, $true, 1, $false, "11+87=98", 2) | Out-Null
$d.Content.Find.Execute("49+31=80", $true, $false, $false, $false, $false, $true, 1, $false, "17+61=78", 2) | Out-Null
$d.Content.Find.Execute("66-34=32", $true, $false, $false, $false, $false, $true, 1, $false, "18+68=86", 2) | Out-Null
$d.Content.Find.Execute("83-22=61", $true, $false, $false, $false, $false, $true, 1, $false, "95-46=49", 2) | Out-Null
$d.Content.Find.Execute("0+63=63", $true, $false, $false, $false, $false, $true, 1, $false, "66-13=53", 2) | Out-Null
$d.Content.Find.Execute("67-7=60", $true, $false, $false, $false, $false, $true, 1, $false, "86-78=8", 2) | Out-Null
$d.Content.Find.Execute("92-3=89", $true, $false, $false, $false, $false, $true, 1, $false, "55-19=36", 2) | Out-Null
$d.Content.Find.Execute("7+74=81", $true, $false, $false, $false, $false, $true, 1, $false, "79-43=36", 2) | Out-Null
$d.Content.Find.Execute("97-32=65", $true, $false, $false, $false, $false, $true, 1, $false, "9+58=67", 2) | Out-Null
$d.Content.Find.Execute("14+68=82", $true, $false, $false, $false, $false, $true, 1, $false, "51+7=58", 2) | Out-Null
$d.Content.Find.Execute("45-14=31", $true, $false, $false, $false, $false, $true, 1, $false, "92-46=46", 2) | Out-Null
$d.Content.Find.Execute("5+13=18", $true, $false, $false, $false, $false, $true, 1, $false, "66-54=12", 2) | Out-Null
$d.Content.Find.Execute("55+6=61", $true, $false, $false, $false, $false, $true, 1, $false, "9+45=54", 2) | Out-Null
$d.Content.Find.Execute("33+25=58", $true, $false, $false, $false, $false, $true, 1, $false, "32+9=41", 2) | Out-Null
$d.Content.Find.Execute("44-4=40", $true, $false, $false, $false, $false, $true, 1, $false, "58-15=43", 2) | Out-Null
$d.Content.Find.Execute("30+63=93", $true, $false, $false, $false, $false, $true, 1, $false, "77-43=34", 2) | Out-Null
$d.Content.Find.Execute("48+9=57", $true, $false, $false, $false, $false, $true, 1, $false, "13+66=79", 2) | Out-Null
$d.Content.Find.Execute("92-22=70", $true, $false, $false, $false, $false, $true, 1, $false, "72+12=84", 2) | Out-Null
$d.Content.Find.Execute("86-33=53", $true, $false, $false, $false, $false, $true, 1, $false, "51-15=36", 2) | Out-Null
$d.Content.Find.Execute("54-24=30", $true, $false, $false, $false, $false, $true, 1, $false, "84-6=78", 2) | Out-Null
$d.Content.Find.Execute("2+49=51", $true, $false, $false, $false, $false, $true, 1, $false, "17+72=89", 2) | Out-Null
$d.Content.Find.Execute("0+6=6", $true, $false, $false, $false, $false, $true, 1, $false, "40-14=26", 2) | Out-Null
